$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

$emuPerPoint = 12700

# Shape "Rectangle 11" (id=12, text "GO"): move up to new y offset
$goShape = $s.Shapes.Item("Rectangle 11")
$goShape.Top = 2684102 / $emuPerPoint

# Shape "Rectangle 24" (id=25, text "Enter %"): move down to new y offset
$pctShape = $s.Shapes.Item("Rectangle 24")
$pctShape.Top = 1828117 / $emuPerPoint

# Shape "Rectangle 32" (id=33, text "Enter $"): remove entirely
$dollarShape = $s.Shapes.Item("Rectangle 32")
$dollarShape.Delete()
